$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1876.8
$ws.Range("I80").Value = 672.3333
$ws.Range("J80").Value = 2679.7778
$ws.Range("K80").Value = 2016.9999
$ws.Range("L80").Value = 8039.3334
$ws.Range("M80").Value = -1018.9999
$ws.Range("N80").Value = -10035.3334
$ws.Range("H83").Value = 1876.8
$ws.Range("I83").Value = 672.3333
$ws.Range("J83").Value = 2679.7778
$ws.Range("K83").Value = 6050.9997
$ws.Range("L83").Value = 24118.0002
$ws.Range("M83").Value = -1058.9997
$ws.Range("N83").Value = -34102.00019999999
$ws.Range("H137").Value = 4081.3103
$ws.Range("I137").Value = 1344.9231
$ws.Range("J137").Value = 27796.666
$ws.Range("K137").Value = 4034.7693
$ws.Range("L137").Value = 83389.99800000001
$ws.Range("M137").Value = -1484.7693
$ws.Range("N137").Value = -88489.99800000001
$ws.Range("H138").Value = 6194.604
$ws.Range("I138").Value = 1204.6522
$ws.Range("J138").Value = 10785.36
$ws.Range("K138").Value = 3613.9566
$ws.Range("L138").Value = 32356.08
$ws.Range("M138").Value = 1526.0434
$ws.Range("N138").Value = -42636.08
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4410.9165
$ws.Range("I2").Value = 3203.2856
$ws.Range("K2").Value = 3203.2856
$ws.Range("M2").Value = -3090.2856
$ws.Range("H76").Value = 108046.25
$ws.Range("J76").Value = 130729
$ws.Range("L76").Value = 130729
$ws.Range("N76").Value = -131405
$ws.Range("H79").Value = 108046.25
$ws.Range("J79").Value = 130729
$ws.Range("L79").Value = 130729
$ws.Range("N79").Value = -133069
$ws.Range("H110").Value = 5761.048
$ws.Range("I110").Value = 5969.8335
$ws.Range("K110").Value = 5969.8335
$ws.Range("M110").Value = -3924.8335
$ws.Range("H116").Value = 4410.9165
$ws.Range("I116").Value = 3203.2856
$ws.Range("K116").Value = 3203.2856
$ws.Range("M116").Value = -909.2856000000002
$ws.Range("H128").Value = 100143
$ws.Range("J128").Value = 100143
$ws.Range("L128").Value = 100143
$ws.Range("N128").Value = -110103
$ws.Range("H132").Value = 21771254
$ws.Range("I132").Value = 1985.8572
$ws.Range("J132").Value = 55634556
$ws.Range("K132").Value = 5957.571599999999
$ws.Range("L132").Value = 166903668
$ws.Range("M132").Value = -3427.571599999999
$ws.Range("N132").Value = -166908728
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4410.9165
$ws.Range("I3").Value = 3203.2856
$ws.Range("K3").Value = 3203.2856
$ws.Range("M3").Value = -3089.2856
$ws.Range("H86").Value = 6508.8213
$ws.Range("I86").Value = 7795.1055
$ws.Range("J86").Value = 3793.3333
$ws.Range("K86").Value = 7795.1055
$ws.Range("L86").Value = 3793.3333
$ws.Range("M86").Value = -6672.1055
$ws.Range("N86").Value = -6039.3333
$ws.Range("H89").Value = 6508.8213
$ws.Range("I89").Value = 7795.1055
$ws.Range("J89").Value = 3793.3333
$ws.Range("K89").Value = 38975.5275
$ws.Range("L89").Value = 18966.6665
$ws.Range("M89").Value = -33359.5275
$ws.Range("N89").Value = -30198.6665
$ws.Range("H105").Value = 5223.1055
$ws.Range("I105").Value = 8149.625
$ws.Range("J105").Value = 3094.7273
$ws.Range("K105").Value = 8149.625
$ws.Range("L105").Value = 3094.7273
$ws.Range("M105").Value = -6402.625
$ws.Range("N105").Value = -6588.7273
$ws.Range("H134").Value = 2812.1667
$ws.Range("I134").Value = 2868.1924
$ws.Range("J134").Value = 2448
$ws.Range("K134").Value = 8604.5772
$ws.Range("L134").Value = 7344
$ws.Range("M134").Value = -6069.5772
$ws.Range("N134").Value = -12414
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 907.5333000000001
$ws.Range("I107").Value = 875.2
$ws.Range("J107").Value = 972.2
$ws.Range("K107").Value = 875.2
$ws.Range("L107").Value = 972.2
$ws.Range("M107").Value = 1044.8
$ws.Range("N107").Value = -4812.2
$ws.Range("H130").Value = 45550
$ws.Range("J130").Value = 45550
$ws.Range("L130").Value = 45550
$ws.Range("N130").Value = -55590
$ws.Range("H131").Value = 25798
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 797.6667
$ws.Range("J68").Value = 797.6667
$ws.Range("L68").Value = 2393.0001
$ws.Range("N68").Value = -4015.0001
$ws.Range("H71").Value = 797.6667
$ws.Range("J71").Value = 797.6667
$ws.Range("L71").Value = 7179.0003
$ws.Range("N71").Value = -15291.0003
$ws.Range("H97").Value = 1324.875
$ws.Range("J97").Value = 1353.4286
$ws.Range("L97").Value = 4060.2858
$ws.Range("N97").Value = -5052.2858
$ws.Range("H129").Value = 1566.4445
$ws.Range("I129").Value = 899.3333
$ws.Range("J129").Value = 1900
$ws.Range("K129").Value = 2697.9999
$ws.Range("L129").Value = 5700
$ws.Range("M129").Value = 2302.0001
$ws.Range("N129").Value = -15700
$ws.Range("H131").Value = 1528.6
$ws.Range("I131").Value = 1349.8
$ws.Range("J131").Value = 1707.4
$ws.Range("K131").Value = 4049.4
$ws.Range("L131").Value = 5122.200000000001
$ws.Range("M131").Value = 990.6000000000004
$ws.Range("N131").Value = -15202.2
$ws.Range("H134").Value = 3788.4285
$ws.Range("J134").Value = 20000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2727.739
$ws.Range("I102").Value = 1421.5625
$ws.Range("J102").Value = 5713.2856
$ws.Range("K102").Value = 1421.5625
$ws.Range("L102").Value = 5713.2856
$ws.Range("M102").Value = 200.4375
$ws.Range("N102").Value = -8957.285599999999
$ws.Range("H122").Value = 83337180
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 100004500
$ws.Range("K122").Value = 1800
$ws.Range("L122").Value = 300013500
$ws.Range("M122").Value = 650
$ws.Range("N122").Value = -300018400
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4515.077
$ws.Range("I7").Value = 4437.2
$ws.Range("K7").Value = 4437.2
$ws.Range("M7").Value = -4325.2
$ws.Range("H22").Value = 2175.6924
$ws.Range("I22").Value = 1968.7778
$ws.Range("J22").Value = 2641.25
$ws.Range("K22").Value = 1968.7778
$ws.Range("L22").Value = 2641.25
$ws.Range("M22").Value = -1673.7778
$ws.Range("N22").Value = -3231.25
$ws.Range("H27").Value = 2175.6924
$ws.Range("I27").Value = 1968.7778
$ws.Range("J27").Value = 2641.25
$ws.Range("K27").Value = 1968.7778
$ws.Range("L27").Value = 2641.25
$ws.Range("M27").Value = -1861.7778
$ws.Range("N27").Value = -2855.25
$ws.Range("H81").Value = 76863.336
$ws.Range("J81").Value = 75545
$ws.Range("L81").Value = 75545
$ws.Range("N81").Value = -77541
$ws.Range("H84").Value = 76863.336
$ws.Range("J84").Value = 75545
$ws.Range("L84").Value = 226635
$ws.Range("N84").Value = -236619
$ws.Range("H126").Value = 4515.077
$ws.Range("I126").Value = 4437.2
$ws.Range("K126").Value = 13311.6
$ws.Range("M126").Value = -10841.6
$ws.Range("H136").Value = 1005401.44
$ws.Range("I136").Value = 1431861.6
$ws.Range("K136").Value = 4295584.800000001
$ws.Range("M136").Value = -4293034.800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 38500
$ws.Range("I56").Value = 41750
$ws.Range("J56").Value = 32000
$ws.Range("K56").Value = 41750
$ws.Range("L56").Value = 32000
$ws.Range("M56").Value = -41036
$ws.Range("N56").Value = -33428
$ws.Range("H138").Value = 115214.5
$ws.Range("J138").Value = 130429
$ws.Range("L138").Value = 130429
$ws.Range("N138").Value = -140709

Write-Output "Applied 195 cell updates across 8 sheets"
